# Update odds values on Sheet1 to reflect refreshed FlashScore data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("G6").Value = 1.85
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 2.63
$ws.Range("K6").Value = 1.83
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("Q6").Value = 2.88
$ws.Range("R6").Value = 1.4
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 2.1
$ws.Range("AA6").Value = 10
$ws.Range("AE6").Value = 5.5
$ws.Range("AJ6").Value = 9

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.63
$ws.Range("L7").Value = 6.5
$ws.Range("Q7").Value = 3.1
$ws.Range("R7").Value = 1.36
$ws.Range("V7").Value = 2
$ws.Range("Z7").Value = 6.5
$ws.Range("AC7").Value = 21
$ws.Range("AG7").Value = 29
$ws.Range("AL7").Value = 21
$ws.Range("AM7").Value = 67

# Row 11
$ws.Range("G11").Value = 1.65
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 5.75
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("AB11").Value = 12
$ws.Range("AE11").Value = 7
$ws.Range("AJ11").Value = 11
$ws.Range("AM11").Value = 67
$ws.Range("AR11").Value = 3.85

# Row 29
$ws.Range("G29").Value = 1.44
$ws.Range("H29").Value = 4.75
$ws.Range("I29").Value = 7
$ws.Range("J29").Value = 1.95
$ws.Range("K29").Value = 2.4
$ws.Range("L29").Value = 6.5
$ws.Range("M29").Value = 1.03
$ws.Range("N29").Value = 15
$ws.Range("O29").Value = 1.2
$ws.Range("P29").Value = 4.33
$ws.Range("Q29").Value = 1.65
$ws.Range("R29").Value = 2.2
$ws.Range("S29").Value = 2.63
$ws.Range("T29").Value = 1.44

# Row 30
$ws.Range("G30").Value = 2.38
$ws.Range("I30").Value = 3
$ws.Range("N30").Value = 9
$ws.Range("Q30").Value = 2.15
$ws.Range("R30").Value = 1.67
$ws.Range("S30").Value = 4
$ws.Range("T30").Value = 1.22
$ws.Range("AB30").Value = 23

# Row 32
$ws.Range("AP32").Value = 1.78
$ws.Range("AQ32").Value = 2.03
